$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.715.45'
$ws.Range("E2").Value = '  -3.05%  '
$ws.Range("D3").Value = '2.097.72'
$ws.Range("E3").Value = '  -1.34%  '
$ws.Range("E4").Value = '  -0.44%  '
$ws.Range("D5").Value = '344.03'
$ws.Range("E5").Value = '  -2.43%  '
$ws.Range("E6").Value = '  -0.33%  '
$ws.Range("D7").Value = '0.5147'
$ws.Range("E7").Value = '  -2.54%  '
$ws.Range("E8").Value = '  -3.14%  '
$ws.Range("D9").Value = '53.07'
$ws.Range("E9").Value = '  -2.08%  '
$ws.Range("D10").Value = '0.09205'
$ws.Range("E10").Value = '  +1.20%  '
$ws.Range("D11").Value = '1.172'
$ws.Range("E11").Value = '  -0.92%  '
$ws.Range("D12").Value = '24.90'
$ws.Range("E12").Value = '  +0.92%  '
$ws.Range("D13").Value = '2.092.24'
$ws.Range("E13").Value = '  -2.07%  '
$ws.Range("D14").Value = '6.762'
$ws.Range("E14").Value = '  -1.59%  '
$ws.Range("D15").Value = '8.191'
$ws.Range("E15").Value = '  +0.74%  '
$ws.Range("D16").Value = '99.48'
$ws.Range("E16").Value = '  -3.01%  '
$ws.Range("E17").Value = '  -2.28%  '
$ws.Range("D18").Value = '1.008'
$ws.Range("E18").Value = '  -0.40%  '
$ws.Range("D19").Value = '20.77'
$ws.Range("E19").Value = '  +6.57%  '
$ws.Range("D20").Value = '0.06637'
$ws.Range("E20").Value = '  -1.16%  '
$ws.Range("D21").Value = '1.006'
$ws.Range("E21").Value = '  -0.42%  '
$ws.Range("D22").Value = '6.189'
$ws.Range("E22").Value = '  -2.67%  '
$ws.Range("D23").Value = '29.758.00'
$ws.Range("E23").Value = '  -3.17%  '
$ws.Range("E24").Value = '  -2.42%  '
$ws.Range("E25").Value = '  -2.89%  '
$ws.Range("D26").Value = '2.339.65'
$ws.Range("E26").Value = '  -1.87%  '
$ws.Range("D27").Value = '21.90'
$ws.Range("E27").Value = '  -2.87%  '
$ws.Range("D28").Value = '162.40'
$ws.Range("E28").Value = '  -1.48%  '
$ws.Range("D29").Value = '2.525'
$ws.Range("E29").Value = '  -1.72%  '
$ws.Range("D30").Value = '132.65'
$ws.Range("E30").Value = '  -2.83%  '
$ws.Range("D31").Value = '1.130'
$ws.Range("E31").Value = '  -5.91%  '
$ws.Range("D32").Value = '0.1049'
$ws.Range("E32").Value = '  -3.24%  '
$ws.Range("E33").Value = '  -1.42%  '
$ws.Range("D34").Value = '6.163'
$ws.Range("E34").Value = '  -3.62%  '
$ws.Range("D35").Value = '3.941'
$ws.Range("E35").Value = '  -1.95%  '
$ws.Range("D36").Value = '6.032'
$ws.Range("E36").Value = '  -2.92%  '
$ws.Range("D37").Value = '10.42'
$ws.Range("D38").Value = '0.02562'
$ws.Range("D39").Value = '0.06714'
$ws.Range("E39").Value = '  -2.62%  '
$ws.Range("D40").Value = '12.44'
$ws.Range("E40").Value = '  -1.02%  '
$ws.Range("D41").Value = '0.2237'
$ws.Range("E41").Value = '  -3.86%  '
$ws.Range("D42").Value = '0.6852'
$ws.Range("E42").Value = '  -1.09%  '
$ws.Range("E43").Value = '  +1.03%  '
$ws.Range("D44").Value = '0.6659'
$ws.Range("E44").Value = '  +3.07%  '
$ws.Range("D45").Value = '14.23'
$ws.Range("E45").Value = '  -3.60%  '
$ws.Range("D46").Value = '2.290'
$ws.Range("E46").Value = '  -2.23%  '
$ws.Range("D47").Value = '3.621'
$ws.Range("E47").Value = '  -3.69%  '
$ws.Range("E48").Value = '  -4.48%  '
$ws.Range("D49").Value = '1.219'
$ws.Range("E49").Value = '  -3.14%  '
$ws.Range("D50").Value = '81.97'
$ws.Range("E50").Value = '  -1.34%  '
$ws.Range("D51").Value = '0.3272'
$ws.Range("E51").Value = '  -3.43%  '
